$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.448.36"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "3.510.96"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.11"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.94"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("D7").Value = "3.508.87"
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.82"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").Value = "4.109.17"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000197"
$ws.Range("E14").Value = "  -5.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.35"
$ws.Range("E15").Value = "  -5.56%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.117"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.489.77"
$ws.Range("E17").Value = "  -3.27%  "
$ws.Range("D18").Value = "65.458.90"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.01"
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.23"
$ws.Range("E21").Value = "  -5.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "414.57"
$ws.Range("E22").Value = "  -4.12%  "
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.32"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").Value = "3.654.13"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -4.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.42"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.67"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.78"
$ws.Range("E30").Value = "  -4.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.01"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "3.513.64"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.21"
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.48"
$ws.Range("E36").Value = "  -4.72%  "
$ws.Range("E37").Value = "  -11.31%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.24"
$ws.Range("E38").Value = "  -6.88%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "172.38"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  -9.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0814"
$ws.Range("E41").Value = "  -4.64%  "
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.854"
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.16"
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("E45").Value = "  -8.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  -7.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.99"
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.84"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.07"
$ws.Range("E50").Value = "  -10.24%  "
$ws.Range("D51").Value = "2.349.28"
$ws.Range("E51").Value = "  -1.21%  "
